# Generate Report for Handback
# Updates the localization-status workbook to reflect that the zh-cn and
# de-de handoffs have been handed back and are in sync with en-US:
#  - Status columns move from "Ready for handoff" to "Handed back: in sync with en-US"
#  - Latest Target File / Latest Handback File / Latest Handback DateTime get filled in
#  - A hyperlink to the source markdown file is added on the new Latest Target File cell
#  - A few columns are widened to comfortably fit the new, longer text

$wb = $excel.ActiveWorkbook

$mdUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/6f225a900fd181add1a52f776d36bac2b4010c14/e2e/eb7f4660-e7f0-479b-970e-ddae80723102.md"
$mdName = "eb7f4660-e7f0-479b-970e-ddae80723102.md"
$statusText = "Handed back: in sync with en-US"

# Excel's ColumnWidth setter rounds to its own internal grid, so request the
# value whose round-tripped width lands as close as possible to the target.
$wideWidth = 29.166666666666668
$maxWidth  = 39.166666666666664

### Overview sheet -----------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Columns.Item(5).ColumnWidth = $wideWidth
$wsOverview.Columns.Item(6).ColumnWidth = $wideWidth
# The zh-cn / de-de status cells on the overview tab mirror the per-language
# "Status" cell, so they flip to the same new text.
$wsOverview.Range("E2").Value = $statusText
$wsOverview.Range("F2").Value = $statusText

### zh-cn sheet ----------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Columns.Item(3).ColumnWidth = $wideWidth
$wsZhCn.Columns.Item(9).ColumnWidth = $maxWidth
$wsZhCn.Columns.Item(10).ColumnWidth = $maxWidth

$wsZhCn.Range("C2").Value = $statusText
$wsZhCn.Range("I2").Value = $mdName
$wsZhCn.Hyperlinks.Add($wsZhCn.Range("I2"), $mdUrl, "", "", $mdName)
$wsZhCn.Range("J2").Value = "eb7f4660-e7f0-479b-970e-ddae80723102.02d2ff92687eb9593844a112278261a058fbddde.zh-cn.xlf"
$wsZhCn.Range("K2").Value = "2016-08-27 08:58:33"

### de-de sheet ----------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Columns.Item(3).ColumnWidth = $wideWidth
$wsDeDe.Columns.Item(9).ColumnWidth = $maxWidth
$wsDeDe.Columns.Item(10).ColumnWidth = $maxWidth

$wsDeDe.Range("C2").Value = $statusText
$wsDeDe.Range("I2").Value = $mdName
$wsDeDe.Hyperlinks.Add($wsDeDe.Range("I2"), $mdUrl, "", "", $mdName)
$wsDeDe.Range("J2").Value = "eb7f4660-e7f0-479b-970e-ddae80723102.02d2ff92687eb9593844a112278261a058fbddde.de-de.xlf"
$wsDeDe.Range("K2").Value = "2016-08-27 08:58:39"
